# Weekly price-log update: a new week's price observation is inserted
# as a new row (221), pushing the existing rows 221-267 down to 222-268.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 221 (shifts 221..267 down to 222..268, extending the
# used range from A1:R267 to A1:R268, and copying row 221's formatting,
# including the date-formatted style on column D).
$ws.Rows("221").Insert()

# Populate the newly inserted row 221 with the new weekly record. Most
# columns mirror the record that is now in row 222 (same market / region /
# category / quality / unit / classification); only the date and the
# volume/price columns change for this week's entry.
$ws.Range("A221").Value = 10
$ws.Range("B221").Value = "Vega Modelo de Temuco"
$ws.Range("C221").Value = "La Araucanía"
$ws.Range("D221").Value = 44476
$ws.Range("E221").Value = 9
$ws.Range("F221").Value = 100112023
$ws.Range("G221").Value = "Brócoli"
$ws.Range("H221").Value = "Sin especificar"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 3000
$ws.Range("K221").Value = 800
$ws.Range("L221").Value = 800
$ws.Range("M221").Value = 800
$ws.Range("N221").Value = "$/unidad"
$ws.Range("O221").Value = "Región Metropolitana"
$ws.Range("P221").Value = 800
$ws.Range("Q221").Value = 1
$ws.Range("R221").Value = "Hortaliza"
